# MOS-23045: Update Master Data as per 22 April Changes
# Append new Postal Code location rows (BNMR children) to master-location.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(10110, 10110, 5, "Postal Code",   "BNMR", "eng", $true, "superadmin", "now()"),
    @(10111, 10111, 5, "Postal Code",   "BNMR", "eng", $true, "superadmin", "now()"),
    @(10113, 10113, 5, "Postal Code",   "BNMR", "eng", $true, "superadmin", "now()"),
    @(10114, 10114, 5, "Postal Code",   "BNMR", "eng", $true, "superadmin", "now()"),
    @(10111, 10111, 5, "code postal",   "BNMR", "fra", $true, "superadmin", "now()"),
    @(10110, 10110, 5, "code postal",   "BNMR", "fra", $true, "superadmin", "now()"),
    @(10113, 10113, 5, "code postal",   "BNMR", "fra", $true, "superadmin", "now()"),
    @(10114, 10114, 5, "code postal",   "BNMR", "fra", $true, "superadmin", "now()"),
    @(10111, 10111, 5, "الرمز البريدي", "BNMR", "ara", $true, "superadmin", "now()"),
    @(10110, 10110, 5, "الرمز البريدي", "BNMR", "ara", $true, "superadmin", "now()")
)

$startRow = 110
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
}

# Clear the stale selection left over from editing (A110:XFD112) so the
# sheet view reverts to a plain default selection, matching the saved file.
$null = $ws.Range("A1").Select()
